# Apply the StructureDefinition metadata refresh:
#  - URL now points at linuxforhealth.org instead of ibm.com
#  - Version bumped 7.0.0 -> 8.0.0
#  - Date advanced to the new publication timestamp
#  - Publisher renamed from "Alvearie Team" to "LinuxForHealth Team"
#  - The stale combined cpt-2/ele-1 constraint text on the root
#    ContactPoint element's Constraint(s) column is cleared

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cdm-communication-contact-point"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
